$d = $word.ActiveDocument

# Replace the first occurrence of $find inside $range with $replace, without
# disturbing neighbouring run formatting or structural anchors
# (commentRangeStart/End, bookmarks, etc.) that sit adjacent to the match.
# We do this by first locating the match (no replace), inserting the new
# text right after the match (so it lands next to whatever follows,
# inheriting neutral/default formatting), and only then deleting the
# original matched text.
function Replace-Text($range, $find, $replace) {
    $r = $range.Duplicate
    $found = $r.Find.Execute($find, $true, $false, $false, $false, $false,
                              $true, 1, $false, "", 0)
    if (-not $found) {
        throw "Text not found: $find"
    }
    $insertPoint = $d.Range($r.End, $r.End)
    $insertPoint.InsertBefore($replace)
    $r.Text = ""
}

# "Thank you for submitting your documents" heading
Replace-Text $d.Paragraphs.Item(14).Range "Thank you for submitting your documents" "Cảm ơn bạn đã gửi các giấy tờ cần thiết"

# "Hi " greeting
Replace-Text $d.Paragraphs.Item(16).Range "Hi " "Xin chào "

# Paragraph about documents / event
Replace-Text $d.Paragraphs.Item(18).Range "Thank you for providing us with your documents for the upcoming " "Cảm ơn bạn đã gửi cho chúng tôi các giấy tờ cần thiết của bạn cho sự kiện "
Replace-Text $d.Paragraphs.Item(18).Range ". Based on the information you’ve given us, we’ll make the necessary arrangements, including accommodation and transportation." " sắp tới. Dựa trên thông tin bạn đã cung cấp, chúng tôi sẽ tiến hành sắp xếp chỗ ở và phương tiện đi lại trong quá trình bạn tham gia sự kiện."

# Currently reviewing documents
Replace-Text $d.Paragraphs.Item(19).Range "We’re currently reviewing your documents and will reach out to you if we need anything else. " "Hiện chúng tôi đang kiểm tra giấy tờ của bạn và sẽ liên hệ với bạn nếu chúng tôi cần thêm thông tin. "

# Contact us via live chat or WhatsApp
Replace-Text $d.Paragraphs.Item(20).Range "If you have any questions, please contact us via " "Nếu bạn cần hỗ trợ, vui lòng liên hệ với chúng tôi qua "
Replace-Text $d.Paragraphs.Item(20).Range " or " " hoặc "

# Contact your country manager paragraph
Replace-Text $d.Paragraphs.Item(21).Range "If you have any questions, please contact your country manager, " "Nếu bạn có bất kỳ thắc mắc nào, vui lòng liên hệ với giám đốc phụ trách quốc gia của bạn "
Replace-Text $d.Paragraphs.Item(21).Range ", at " ", qua email "
Replace-Text $d.Paragraphs.Item(21).Range " or " " hoặc số "

# Look forward to seeing you
Replace-Text $d.Paragraphs.Item(22).Range "We look forward to seeing you at " "Chúng tôi rất mong được gặp bạn tại sự kiện "
